$d = $word.ActiveDocument
$p32 = $d.Paragraphs(32)
$r = $d.Range($p32.Range.End - 1, $p32.Range.End - 1)
$r.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>  <w:pPr>    <w:spacing w:line="400" w:lineRule="exact"/>    <w:ind w:firstLine="420"/>    <w:rPr>      <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia"/>      <w:sz w:val="24"/>      <w:szCs w:val="24"/>    </w:rPr>  </w:pPr>  <w:r>    <w:rPr>      <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:hint="eastAsia"/>      <w:sz w:val="24"/>      <w:szCs w:val="24"/>    </w:rPr>    <w:lastRenderedPageBreak/>    <w:t>最后</w:t>  </w:r>  <w:r>    <w:rPr>      <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:hint="eastAsia"/>      <w:sz w:val="24"/>      <w:szCs w:val="24"/>    </w:rPr>    <w:t>，在</w:t>  </w:r>  <w:r>    <w:rPr>      <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:hint="eastAsia"/>      <w:sz w:val="24"/>      <w:szCs w:val="24"/>    </w:rPr>    <w:t>数据存储与可视化</w:t>  </w:r>  <w:r>    <w:rPr>      <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:hint="eastAsia"/>      <w:sz w:val="24"/>      <w:szCs w:val="24"/>    </w:rPr>    <w:t>分析方面，研究了项目使用的本地存储技术，包括localStorage存储和sessionStorage；组件间数据传输方式，包括</w:t>  </w:r>  <w:r>    <w:rPr>      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>      <w:sz w:val="24"/>      <w:szCs w:val="24"/>    </w:rPr>    <w:t>react</w:t>  </w:r>  <w:r>    <w:rPr>      <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:hint="eastAsia"/>      <w:sz w:val="24"/>      <w:szCs w:val="24"/>    </w:rPr>    <w:t>特有的传递数据方式时</w:t>  </w:r>  <w:r>    <w:rPr>      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>      <w:sz w:val="24"/>      <w:szCs w:val="24"/>    </w:rPr>    <w:t>HTML5</w:t>  </w:r>  <w:r>    <w:rPr>      <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:hint="eastAsia"/>      <w:sz w:val="24"/>      <w:szCs w:val="24"/>    </w:rPr>    <w:t>下的传递方式；数据可视化分析。</w:t>  </w:r>  <w:bookmarkStart w:id="0" w:name="_GoBack"/>  <w:bookmarkEnd w:id="0"/></w:p><w:p>  <w:pPr>    <w:spacing w:line="400" w:lineRule="exact"/>    <w:ind w:firstLine="420"/>    <w:rPr>      <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:hint="eastAsia"/>      <w:sz w:val="24"/>      <w:szCs w:val="24"/>    </w:rPr>  </w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
Write-Host "AFTER count:" $d.Paragraphs.Count
for ($i = 30; $i -le 37; $i++) {
  Write-Host $i $d.Paragraphs($i).Range.Text.Length
}
